# Add "SD Card" line item to the budget sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row: A8 = "SD Card", B8 = 40 (cost, currency-formatted style already present)
$ws.Range("A8").Value = "SD Card"
$ws.Range("B8").Value = 40

# Update the active selection to reflect A9 (next empty row), matching the diff.
$ws.Range("A9").Select()
